# Insert a new price-observation row for "Frambuesa" (Raspberry) at row 138,
# shifting all the existing rows 138..227 down to 139..228, and populate
# the newly-opened row 138 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlShiftDown = -4121
$ws.Rows("138:138").Insert(-4121)

$ws.Range("A138").Value = 6
$ws.Range("B138").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C138").Value = "Metropolitana"
$ws.Range("D138").Value = 44907
$ws.Range("E138").Value = 13
$ws.Range("F138").Value = "Fruta"
$ws.Range("G138").Value = 100101
$ws.Range("H138").Value = "Berries"
$ws.Range("I138").Value = 100101004
$ws.Range("J138").Value = "Frambuesa"
$ws.Range("K138").Value = "Sin especificar"
$ws.Range("L138").Value = "Especial"
$ws.Range("M138").Value = 500
$ws.Range("N138").Value = 8000
$ws.Range("O138").Value = 8000
$ws.Range("P138").Value = 8000
$ws.Range("Q138").Value = "`$/bandeja 2 kilos"
$ws.Range("R138").Value = "Región del Maule"
$ws.Range("S138").Value = 4000
$ws.Range("T138").Value = 2
